# Add a new attendance sheet for 2020-11-14, modelled on the most recent
# existing day sheet (2020-11-12), and fill it with that day's readings.

$wb = $excel.ActiveWorkbook

# The previous day's sheet is the template: same headers, same per-cell
# styling (bold/centered/bordered header row + "Sr. No" column).
$template = $wb.Worksheets.Item("2020-11-12")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Duplicate it and park the copy after the last existing tab, then rename.
$template.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "2020-11-14"

# The template sheet spans A1:G6 (header + 5 data rows); the new day only
# has header + 4 data rows (A1:G5), so drop the now-unused last row
# entirely (value + formatting) rather than just clearing its contents.
$ws.Rows.Item(6).Delete()

# Overwrite with the new day's readings (values only -- formatting/styles
# already carried over from the template sheet for row 1 / column A).
$headers = @("Sr. No", "Name", "Address", "Job", "Time-Stamp", "SpO2_value", "Heart-rate")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$data = @(
    @(1, "sachin", "xyz/xyz", "coder", "16:07:22", 97.75047190698253, 0),
    @(1, "sachin", "xyz/xyz", "coder", "16:09:49", 94.55808920791817, 57.20660651222352),
    @(1, "sachin", "xyz/xyz", "coder", "16:10:54", 93.31806497949607, 69.20471482450479),
    @(1, "sachin", "xyz/xyz", "coder", "16:12:47", 94.42039855332457, 51.60096886331338)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

Write-Host "Added sheet '2020-11-14' with $($data.Length) attendance rows."
